# LITE-28847: Support 'late' field in product capabilities
#
# Inserts a new "Pay-as-you-go late charges support" capability row into the
# "Capabilities" sheet (row 5), pushing the existing rows down by one and
# leaving the sheet's built-in list data validations to extend over the
# newly shifted ranges automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capabilities")

# Insert a new blank row above the current row 5 ("Consumption reporting for
# Reservation Items"), shifting it (and everything below) down one row.
$ws.Rows.Item(5).Insert()

# Populate the new row with the new capability.
$ws.Cells.Item(5, 1).Value = "Pay-as-you-go late charges support"
$ws.Cells.Item(5, 2).Value = "-"
$ws.Cells.Item(5, 3).Value = "Disabled"

# Match the author's final selection (cell A10 on the Capabilities sheet).
$ws.Activate()
$null = $ws.Range("A10").Select()
